$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new LCSC product-detail URL to G8 (previously empty) and
# enable text wrapping on it.
$ws.Range("G8").Value = "https://lcsc.com/product-detail/USB-Connectors_Jing-Extension-of-the-Electronic-Co-920-E52A2021S10100_C10418.html"
$ws.Range("G8").WrapText = $true

# The existing hyperlink cell (G6, pogo-pin supplier link) also gets
# wrapped now that the column is wider.
$ws.Range("G6").WrapText = $true

# Widen column G so the long URLs wrap onto a few lines instead of one.
$ws.Columns("G").ColumnWidth = 36.3

# Let the affected rows grow to fit the wrapped text.
$ws.Rows(6).RowHeight = 46
$ws.Rows(8).RowHeight = 53

# Re-select cell E20 (matches the saved view state of the edited file).
$ws.Range("E20").Select()
